# Utah overview workbook: convert numeric "count" and percent/currency
# cells that were stored as numbers into literal text cells (matching
# COMM's text-edit pass), and append a new "Total" row to the County sheet.

function Set-TextValue($range, $value) {
    # Force Excel to store the value as literal text (inlineStr/shared
    # string) rather than re-parsing it back into a number/currency/percent,
    # then drop back to the Normal style so no stray number-format style
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overall sheet
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Range("A2") "718"

# ---------------------------------------------------------------------
# County sheet
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

Set-TextValue $wsCounty.Range("B2") "6"
Set-TextValue $wsCounty.Range("B3") "31"
Set-TextValue $wsCounty.Range("B4") "8"
Set-TextValue $wsCounty.Range("B5") "38"
Set-TextValue $wsCounty.Range("B6") "2"
Set-TextValue $wsCounty.Range("B7") "4"
Set-TextValue $wsCounty.Range("B8") "1"
Set-TextValue $wsCounty.Range("B9") "18"
Set-TextValue $wsCounty.Range("B10") "8"
Set-TextValue $wsCounty.Range("B11") "1"
Set-TextValue $wsCounty.Range("B12") "1"
Set-TextValue $wsCounty.Range("B13") "358"
Set-TextValue $wsCounty.Range("B14") "5"
Set-TextValue $wsCounty.Range("B15") "6"
Set-TextValue $wsCounty.Range("B16") "5"
Set-TextValue $wsCounty.Range("B17") "40"
Set-TextValue $wsCounty.Range("B18") "4"
Set-TextValue $wsCounty.Range("B19") "2"
Set-TextValue $wsCounty.Range("B20") "83"
Set-TextValue $wsCounty.Range("B21") "6"
Set-TextValue $wsCounty.Range("B22") "33"
Set-TextValue $wsCounty.Range("B23") "3"
Set-TextValue $wsCounty.Range("B24") "55"

# Rows 25-28 (Rich, Millard, Piute, Daggett counties) had all-zero stats;
# COMM's edit replaces them with formatted zero text.
$zeroRows = 25, 26, 27, 28
foreach ($row in $zeroRows) {
    Set-TextValue $wsCounty.Range("B$row") "0.00%"
    Set-TextValue $wsCounty.Range("C$row") "$0"
    Set-TextValue $wsCounty.Range("D$row") "0.00%"
    Set-TextValue $wsCounty.Range("E$row") "0.00%"
    Set-TextValue $wsCounty.Range("F$row") "0.00%"
}

# New row 29: county-sheet Total row.
Set-TextValue $wsCounty.Range("A29") "Total"
Set-TextValue $wsCounty.Range("B29") "718"
Set-TextValue $wsCounty.Range("C29") "$1,665,786,876"
Set-TextValue $wsCounty.Range("D29") "9.73%"
Set-TextValue $wsCounty.Range("E29") "-23.91%"
Set-TextValue $wsCounty.Range("F29") "72.14%"

# ---------------------------------------------------------------------
# Congressional District sheet
# ---------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
Set-TextValue $wsCd.Range("B2") "206"
Set-TextValue $wsCd.Range("B3") "224"
Set-TextValue $wsCd.Range("B4") "191"
Set-TextValue $wsCd.Range("B5") "97"
Set-TextValue $wsCd.Range("B6") "718"

# ---------------------------------------------------------------------
# Size sheet
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
Set-TextValue $wsSize.Range("B2") "183"
Set-TextValue $wsSize.Range("B3") "242"
Set-TextValue $wsSize.Range("B4") "99"
Set-TextValue $wsSize.Range("B5") "82"
Set-TextValue $wsSize.Range("B6") "80"
Set-TextValue $wsSize.Range("B7") "32"
Set-TextValue $wsSize.Range("B8") "718"

# ---------------------------------------------------------------------
# Subsector sheet
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
Set-TextValue $wsSub.Range("B2") "93"
Set-TextValue $wsSub.Range("B3") "159"
Set-TextValue $wsSub.Range("B4") "45"
Set-TextValue $wsSub.Range("B5") "67"
Set-TextValue $wsSub.Range("B6") "5"
Set-TextValue $wsSub.Range("B7") "177"
Set-TextValue $wsSub.Range("B8") "7"
Set-TextValue $wsSub.Range("B9") "39"
Set-TextValue $wsSub.Range("B10") "4"
Set-TextValue $wsSub.Range("B11") "119"
Set-TextValue $wsSub.Range("B12") "3"
Set-TextValue $wsSub.Range("B13") "718"
